# "aggiornamento fino a 21 marzo" - append 4 new daily rows (230-233) to
# the Spilamberto report sheet, extending the dimension to A1:D233.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$newRows = @(
    @(44304, 2, 23, 180.6046329014527),
    @(44305, 4, 20, 157.0475068708284),
    @(44306, 7, 27, 212.0141342756183),
    @(44307, 0, 27, 212.0141342756183)
)

$startRow = 230
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $data[0]

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
